# Cypress Ascendant Services LLC sample template - Ops update
# 12/30/2024 new changes in ops (ordercreation & orderpage & order form)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Cell value updates (OrderID / Emp ID / Typist / Typist QC)
# ---------------------------------------------------------------
$ws.Range("B2").Value = "001CAS"
$ws.Range("B3").Value = "002CAS"

$ws.Range("C3").Value = "SIPL5316"

$ws.Range("E2").Value = "SIPL0102"
$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F2").Value = "SIPL0103"
$ws.Range("F3").Value = "SIPL0103"

# ---------------------------------------------------------------
# 2. Formatting: unify the data rows (2:3) font color to solid
#    black and give the date column the built-in date/time format.
# ---------------------------------------------------------------
$dataRange = $ws.Range("A2:M3")
$dataRange.Font.Color = 0

$ws.Range("A2:A3").NumberFormat = "m/d/yy h:mm"

# ---------------------------------------------------------------
# 3. Header row: bold black text on a gold/tan highlight fill
# ---------------------------------------------------------------
$headerRange = $ws.Range("A1:M1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 0
$headerRange.Interior.Color = 10086143
$headerRange.Interior.PatternColor = 0

# ---------------------------------------------------------------
# 4. Column widths (cosmetic resize following the data edits)
# ---------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 35.5
$ws.Columns.Item(5).ColumnWidth = 11.83
$ws.Columns.Item(6).ColumnWidth = 15.17
$ws.Columns.Item(7).ColumnWidth = 15.17
$ws.Columns.Item(10).ColumnWidth = 15.17
$ws.Columns.Item(12).ColumnWidth = 10

# ---------------------------------------------------------------
# 5. Misc view bookkeeping
# ---------------------------------------------------------------
$ws.Range("E6").Select()
